# Insert one new weekly record for "Macroferia Regional de Talca" / Coliflor.
# The new observation is inserted as row 123, pushing the existing rows
# 123-136 down to 124-137 (dimension grows from A1:R136 to A1:R137).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(123).Insert()

$ws.Range("A123").Value = 5
$ws.Range("B123").Value = "Macroferia Regional de Talca"
$ws.Range("C123").Value = "Maule"
$ws.Range("D123").Value = 44449
$ws.Range("E123").Value = 7
$ws.Range("F123").Value = 100112008
$ws.Range("G123").Value = "Coliflor"
$ws.Range("H123").Value = "Sin especificar"
$ws.Range("I123").Value = "Primera"
$ws.Range("J123").Value = 3000
$ws.Range("K123").Value = 600
$ws.Range("L123").Value = 600
$ws.Range("M123").Value = 600
$ws.Range("N123").Value = "$/unidad"
$ws.Range("O123").Value = "Región del Maule"
$ws.Range("P123").Value = 600
$ws.Range("Q123").Value = 1
$ws.Range("R123").Value = "Hortaliza"
